$p = $ppt.ActivePresentation

# --- Slide 1: "Subtitle 2" placeholder (Team Ninja / Supervised by ...) ---
$s1 = $p.Slides.Item(1)
$subTitle = $s1.Shapes.Item(7)
$tr = $subTitle.TextFrame.TextRange

# Route the paragraph-2 text change through an unrelated placeholder value
# (no shared prefix/suffix with either the old or new text) so the host
# replaces the whole paragraph with a single clean run instead of doing a
# minimal-diff split across runs.
$tr.Paragraphs(2).Text = "XyzPlaceholderXyz"
$tr.Paragraphs(2).Text = "Supervised by Professor Pei Young"

# --- Slide 11: Title "Tools used for the team" (merge two runs into one) ---
$s11 = $p.Slides.Item(11)
$title = $s11.Shapes.Item(1)
$titleTr = $title.TextFrame.TextRange

$titleTr.Text = "XyzPlaceholderXyz"
$titleTr.Text = "Tools used for the team"
